# The document has two logos that are placed as inline pictures in the
# page headers (the BTEC logo, "BTec_Logo-Orange") and page footers (the
# Pearson logo, "...PearsonLogo.png"). Each logo appears twice because the
# section has separate odd/first-page header and footer parts.
#
# The edit renames the pictures:
#   - BTEC logo in both headers:  image2.jpg -> image1.jpg
#   - Pearson logo in both footers: image1.png -> image2.png
#
# The picture's description (alt text) is left untouched; only its Name
# changes. We reach each inline picture through
# Section.Headers/Footers(n).Range.InlineShapes, and we go through the
# Selection object (select the picture's range, then use
# Selection.InlineShapes) because that is the reliable way to reach every
# header/footer picture and set its Name.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-HeaderFooterPicture($headerFooter, $newName) {
    if ($headerFooter.Exists -and $headerFooter.Range.InlineShapes.Count -gt 0) {
        $headerFooter.Range.InlineShapes.Item(1).Range.Select()
        $word.Selection.InlineShapes.Item(1).Name = $newName
    }
}

for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
    Rename-HeaderFooterPicture $sec.Headers.Item($hi) "image1.jpg"
}

for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
    Rename-HeaderFooterPicture $sec.Footers.Item($fi) "image2.png"
}
